$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that currently sits right after
#    the title paragraph ("Play Age of the Gods Epic Troy Slot for Free
#    2021" / "Meta description: Discover the immersive world ..."). We
#    locate it by its distinctive leading text rather than assuming a fixed
#    paragraph index, then delete the whole paragraph (Range.Delete removes
#    the paragraph mark too, so the following paragraph shifts up cleanly).
# ---------------------------------------------------------------------------
$metaPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Meta description")) {
        $metaPara = $p
        break
    }
}
if ($metaPara -ne $null) {
    $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) At the very end of the document there is a paragraph (italic) that used
#    to read "Create a feature image for ...". We need to:
#      a) insert a brand-new bold paragraph right before it that repeats the
#         page title ("Play Age of the Gods Epic Troy Slot for Free 2021"),
#      b) turn the old italic paragraph's text into the meta-description
#         sentence ("Discover the immersive world ...").
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)

# Insert a fresh paragraph after the paragraph that precedes the last one,
# so the new paragraph mark inherits that (plain) paragraph's run
# formatting instead of the italic formatting used by the last paragraph.
$anchorPara = $d.Paragraphs($d.Paragraphs.Count - 1)
$anchorPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs($d.Paragraphs.Count - 1)

# Strip away anything the new paragraph mark may have inherited from the
# anchor paragraph (e.g. the "What we don't like" bullet-list style/indent)
# so it ends up as a plain body paragraph.
$newPara.Range.ParagraphFormat.LeftIndent = 0
$newPara.Style = "Normal"

$newRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$newRange.Text = "Play Age of the Gods Epic Troy Slot for Free 2021"
$newRange.Font.Bold = $true

# Now update the final paragraph's text (still italic) to the new sentence.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$lastRange.Text = "Discover the immersive world of Age of the Gods Epic Troy slot - play now for free and experience high-value combinations and engaging bonus modes."
